# The author's edit re-styles the three tables in the deck (on slides 14,
# 15 and 16) to use the table style {9DFDF103-63D8-4BF3-8B84-76CF84215F48}
# instead of the previous {E97CB087-F1F3-45AA-BB33-6A434CA24C32}.
#
# Each of these three slides contains exactly one shape, and that shape is
# the graphicFrame hosting the <a:tbl>. PowerPoint's object model exposes
# the table's style only through Table.ApplyStyle(styleId) (Table.Style is
# read-only), so that is what we call for each of them.

$p = $ppt.ActivePresentation

$oldStyleId = "{E97CB087-F1F3-45AA-BB33-6A434CA24C32}"
$newStyleId = "{9DFDF103-63D8-4BF3-8B84-76CF84215F48}"

$tableSlideIndexes = @(14, 15, 16)

foreach ($slideIndex in $tableSlideIndexes) {
    $slide = $p.Slides.Item($slideIndex)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $table = $shape.Table
            if ($table.Style -eq $oldStyleId) {
                $table.ApplyStyle($newStyleId)
            }
        }
    }
}
